$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row: Id=3002, Des=背包已满 ("Backpack is full")
$ws.Range("A25").Value = 3002
$ws.Range("B25").Value = "背包已满"

# Grow the table / ListObject so it covers the new row (this also updates
# the worksheet dimension and the table's autoFilter range).
$table = $ws.ListObjects.Item(1)
$table.Resize($ws.Range("A1:B25"))

# Move the active selection to follow the newly-added row, matching the
# author's editing position after inserting the row.
$ws.Range("B23").Select()
